$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1564.579
$ws.Range("J41").Value = 2344.4167
$ws.Range("L41").Value = 2344.4167
$ws.Range("N41").Value = -3224.4167
$ws.Range("H92").Value = 2594.8823
$ws.Range("I92").Value = 199.63637
$ws.Range("J92").Value = 6986.1665
$ws.Range("K92").Value = 199.63637
$ws.Range("L92").Value = 6986.1665
$ws.Range("M92").Value = 1048.36363
$ws.Range("N92").Value = -9482.166499999999
$ws.Range("H100").Value = 5761.8
$ws.Range("I100").Value = 6070.769
$ws.Range("K100").Value = 6070.769
$ws.Range("M100").Value = -5529.769
$ws.Range("H112").Value = 11726.154
$ws.Range("I112").Value = 1000
$ws.Range("K112").Value = 3000
$ws.Range("M112").Value = -1892
$ws.Range("H138").Value = 28481.764
$ws.Range("I138").Value = 1700.04
$ws.Range("K138").Value = 5100.12
$ws.Range("M138").Value = 39.88000000000011
$ws.Range("H141").Value = 7549.5
$ws.Range("I141").Value = 7549.5
$ws.Range("K141").Value = 22648.5
$ws.Range("M141").Value = -17468.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 12511
$ws.Range("I61").Value = 1577.5
$ws.Range("J61").Value = 19800
$ws.Range("K61").Value = 1577.5
$ws.Range("L61").Value = 19800
$ws.Range("M61").Value = -1365.5
$ws.Range("N61").Value = -20224
$ws.Range("H88").Value = 6085
$ws.Range("I88").Value = 1605.3334
$ws.Range("J88").Value = 11460.6
$ws.Range("K88").Value = 1605.3334
$ws.Range("L88").Value = 11460.6
$ws.Range("M88").Value = -1199.3334
$ws.Range("N88").Value = -12272.6
$ws.Range("H91").Value = 6085
$ws.Range("I91").Value = 1605.3334
$ws.Range("J91").Value = 11460.6
$ws.Range("K91").Value = 1605.3334
$ws.Range("L91").Value = 11460.6
$ws.Range("M91").Value = -201.3334
$ws.Range("N91").Value = -14268.6
$ws.Range("H97").Value = 898.56525
$ws.Range("I97").Value = 916.55554
$ws.Range("K97").Value = 916.55554
$ws.Range("M97").Value = -420.55554
$ws.Range("H132").Value = 1902.3043
$ws.Range("I132").Value = 1750.1428
$ws.Range("K132").Value = 5250.428400000001
$ws.Range("M132").Value = -2720.428400000001
$ws.Range("H136").Value = 12511
$ws.Range("I136").Value = 1577.5
$ws.Range("J136").Value = 19800
$ws.Range("K136").Value = 4732.5
$ws.Range("L136").Value = 59400
$ws.Range("M136").Value = -2182.5
$ws.Range("N136").Value = -64500

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1404.8
$ws.Range("I64").Value = 2266.5
$ws.Range("K64").Value = 2266.5
$ws.Range("M64").Value = -2041.5
$ws.Range("H67").Value = 1404.8
$ws.Range("I67").Value = 2266.5
$ws.Range("K67").Value = 2266.5
$ws.Range("M67").Value = -1486.5
$ws.Range("H80").Value = 1017.1539
$ws.Range("I80").Value = 1257.8334
$ws.Range("J80").Value = 810.8570999999999
$ws.Range("K80").Value = 1257.8334
$ws.Range("L80").Value = 810.8570999999999
$ws.Range("M80").Value = -259.8334
$ws.Range("N80").Value = -2806.8571
$ws.Range("H83").Value = 1017.1539
$ws.Range("I83").Value = 1257.8334
$ws.Range("J83").Value = 810.8570999999999
$ws.Range("K83").Value = 6289.166999999999
$ws.Range("L83").Value = 4054.2855
$ws.Range("M83").Value = -1297.166999999999
$ws.Range("N83").Value = -14038.2855
$ws.Range("H86").Value = 1187
$ws.Range("I86").Value = 881.8182
$ws.Range("K86").Value = 881.8182
$ws.Range("M86").Value = 241.1818
$ws.Range("H89").Value = 1187
$ws.Range("I89").Value = 881.8182
$ws.Range("K89").Value = 4409.091
$ws.Range("M89").Value = 1206.909
$ws.Range("H94").Value = 1994.0714
$ws.Range("I94").Value = 916.36365
$ws.Range("J94").Value = 5945.6665
$ws.Range("K94").Value = 916.36365
$ws.Range("L94").Value = 5945.6665
$ws.Range("M94").Value = -465.36365
$ws.Range("N94").Value = -6847.6665
$ws.Range("H99").Value = 1431.5
$ws.Range("I99").Value = 1217.8
$ws.Range("J99").Value = 2500
$ws.Range("K99").Value = 1217.8
$ws.Range("L99").Value = 2500
$ws.Range("M99").Value = 280.2
$ws.Range("N99").Value = -5496
$ws.Range("H105").Value = 3250.2727
$ws.Range("I105").Value = 2108.8
$ws.Range("K105").Value = 2108.8
$ws.Range("M105").Value = -361.8000000000002
$ws.Range("H134").Value = 3624.3125
$ws.Range("I134").Value = 3097.818
$ws.Range("K134").Value = 9293.454000000002
$ws.Range("M134").Value = -6758.454000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4693.125
$ws.Range("I99").Value = 2498.75
$ws.Range("J99").Value = 6887.5
$ws.Range("K99").Value = 2498.75
$ws.Range("L99").Value = 6887.5
$ws.Range("M99").Value = -1000.75
$ws.Range("N99").Value = -9883.5
$ws.Range("H126").Value = 4693.125
$ws.Range("I126").Value = 2498.75
$ws.Range("J126").Value = 6887.5
$ws.Range("K126").Value = 7496.25
$ws.Range("L126").Value = 20662.5
$ws.Range("M126").Value = -5026.25
$ws.Range("N126").Value = -25602.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 5184.905
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 5184.905
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 15554.715
$ws.Range("M80").Value = ""
$ws.Range("N80").Value = -17426.715
$ws.Range("H83").Value = 5184.905
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 5184.905
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 46664.145
$ws.Range("M83").Value = ""
$ws.Range("N83").Value = -56024.145

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3404
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 3404
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 3404
$ws.Range("M113").Value = ""
$ws.Range("N113").Value = -7744
$ws.Range("H122").Value = 3958.3333
$ws.Range("I122").Value = 3687.5
$ws.Range("J122").Value = 4500
$ws.Range("K122").Value = 11062.5
$ws.Range("L122").Value = 13500
$ws.Range("M122").Value = -8612.5
$ws.Range("N122").Value = -18400
$ws.Range("H126").Value = 2000
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -3530

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1619.2941
$ws.Range("I16").Value = 1361.8667
$ws.Range("K16").Value = 1361.8667
$ws.Range("M16").Value = -1191.8667
$ws.Range("H22").Value = 1621.6428
$ws.Range("I22").Value = 1412.875
$ws.Range("K22").Value = 1412.875
$ws.Range("M22").Value = -1117.875
$ws.Range("H27").Value = 1621.6428
$ws.Range("I27").Value = 1412.875
$ws.Range("K27").Value = 1412.875
$ws.Range("M27").Value = -1305.875
$ws.Range("H40").Value = 3116.5
$ws.Range("I40").Value = 2849.8
$ws.Range("K40").Value = 2849.8
$ws.Range("M40").Value = -2713.8
$ws.Range("H61").Value = 3850
$ws.Range("I61").Value = 2750
$ws.Range("K61").Value = 2750
$ws.Range("M61").Value = -2548
$ws.Range("H93").Value = 2240.2727
$ws.Range("I93").Value = 1705.25
$ws.Range("K93").Value = 1705.25
$ws.Range("M93").Value = -457.25
$ws.Range("H113").Value = 3850
$ws.Range("I113").Value = 2750
$ws.Range("K113").Value = 2750
$ws.Range("M113").Value = -580
$ws.Range("H132").Value = 3565.1
$ws.Range("I132").Value = 3123.2917
$ws.Range("J132").Value = 5332.3335
$ws.Range("K132").Value = 9369.875100000001
$ws.Range("L132").Value = 15997.0005
$ws.Range("M132").Value = -6839.875100000001
$ws.Range("N132").Value = -21057.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4444
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = ""
$ws.Range("H65").Value = 4444
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = ""
$ws.Range("H74").Value = 16099.75
$ws.Range("I74").Value = 10569
$ws.Range("J74").Value = 16889.857
$ws.Range("K74").Value = 10569
$ws.Range("L74").Value = 16889.857
$ws.Range("M74").Value = -9633
$ws.Range("N74").Value = -18761.857
$ws.Range("H77").Value = 16099.75
$ws.Range("I77").Value = 10569
$ws.Range("J77").Value = 16889.857
$ws.Range("K77").Value = 31707
$ws.Range("L77").Value = 50669.571
$ws.Range("M77").Value = -27027
$ws.Range("N77").Value = -60029.571
$ws.Range("H105").Value = 28000
$ws.Range("J105").Value = 28000
$ws.Range("L105").Value = 28000
$ws.Range("N105").Value = -34988
$ws.Range("H107").Value = 818.41174
$ws.Range("I107").Value = 774.4
$ws.Range("J107").Value = 1148.5
$ws.Range("K107").Value = 2323.2
$ws.Range("L107").Value = 3445.5
$ws.Range("M107").Value = -403.1999999999998
$ws.Range("N107").Value = -7285.5
$ws.Range("H122").Value = 42848.605
$ws.Range("I122").Value = 45382.387
$ws.Range("K122").Value = 136147.161
$ws.Range("M122").Value = -133697.161
$ws.Range("H132").Value = 1658.02
$ws.Range("I132").Value = 1174.8422
$ws.Range("K132").Value = 3524.5266
$ws.Range("M132").Value = -994.5266000000001
